# aggiornamento fino a 21 marzo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(230, 44304, 3, 11, 111.5392415331576),
    @(231, 44305, 1, 11, 111.5392415331576),
    @(232, 44306, 5, 16, 162.2388967755019),
    @(233, 44307, 1, 17, 172.3788278239708)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item(229, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}
